$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 38
$ws.Range("I2").Value = 122
$ws.Range("J2").Value = 506
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 134
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = 87
$ws.Range("P2").Value = 7
$ws.Range("Q2").Value = 1
$ws.Range("S2").Value = 61
$ws.Range("T2").Value = 91
$ws.Range("U2").Value = 3
$ws.Range("V2").Value = 738
$ws.Range("Z2").Value = 11
$ws.Range("AA2").Value = 5
